# Added fixed effect models
#
# The original single-sheet workbook ("model structure") is split into two
# sheets: the existing intercept-style table (renamed) and a brand new
# "model structure - fixed" sheet describing the new brm_12_het_fixed1-4
# models (treatment main effect + treatment*slope interactions).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# The original sheet keeps all of its data - it only gets a new, more
# specific name now that there is more than one "model structure" sheet.
$ws1.Name = "model structure - Intercept"

# The old "last edited cell" selection (A21) is stale; the author's new
# selection on this sheet is the column-header row.
$ws1.Range("A3:D3").Select()

# New sheet, placed immediately after "model structure - Intercept", and
# becomes the active tab (matches Excel's normal "insert sheet" behaviour).
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "model structure - fixed"

# ---- Header row (row 1) -----------------------------------------------
$ws2.Range("A1").Value = "Model Name"
$ws2.Range("B1").Value = "Intercept"
$ws2.Range("C1").Value = "Trt"
$ws2.Range("D1").Value = "Linear Slope"
$ws2.Range("E1").Value = "Quadratic Slope"
$ws2.Range("F1").Value = "Trt*Linear Slope"
$ws2.Range("G1").Value = "Trt*Quadratic Slope"

# Reuse the existing "big header" cell format (bold, 14pt, centered) from
# the sibling sheet instead of building it up property-by-property, so we
# don't leave unused styles behind in styles.xml.
$ws1.Range("A3").Copy()
foreach ($addr in @("A1", "B1", "D1", "E1", "F1", "G1")) {
    $ws2.Range($addr).PasteSpecial(-4122)
}

# "Trt" gets the smaller bold header style used elsewhere on the original
# sheet (e.g. the "NO MAT" / "NO PE" annotation cells).
$ws1.Range("F11").Copy()
$ws2.Range("C1").PasteSpecial(-4122)

$ws2.Rows.Item(1).RowHeight = 19

# ---- Data rows (centered "X" marks, same style as the original table) --
$ws2.Range("A2").Value = "brm_12_het_fixed1"
$ws2.Range("B2").Value = "X"
$ws2.Range("C2").Value = "X"
$ws2.Range("D2").Value = "X"
$ws2.Range("E2").Value = "X"
$ws2.Range("F2").Value = "X"
$ws2.Range("G2").Value = "X"

$ws2.Range("A3").Value = "brm_12_het_fixed2"
$ws2.Range("B3").Value = "X"
$ws2.Range("C3").Value = "X"
$ws2.Range("D3").Value = "X"
$ws2.Range("E3").Value = "X"
$ws2.Range("F3").Value = "X"

$ws2.Range("A4").Value = "brm_12_het_fixed3"
$ws2.Range("B4").Value = "X"
$ws2.Range("C4").Value = "X"
$ws2.Range("D4").Value = "X"
$ws2.Range("E4").Value = "X"
$ws2.Range("G4").Value = "X"

$ws2.Range("A5").Value = "brm_12_het_fixed4"
$ws2.Range("B5").Value = "X"
$ws2.Range("C5").Value = "X"
$ws2.Range("D5").Value = "X"
$ws2.Range("E5").Value = "X"

# Center-align every cell in the grid (B:G, rows 2-5), including the
# blank ones (G3, F4, F5, G5) - matches the "NO"-style gaps on the
# original sheet where the cell is blank but still centered.
$ws1.Range("B4").Copy()
$ws2.Range("B2:G2").PasteSpecial(-4122)
$ws2.Range("B3:G3").PasteSpecial(-4122)
$ws2.Range("B4:G4").PasteSpecial(-4122)
$ws2.Range("B5:G5").PasteSpecial(-4122)

$ws2.Range("A7").Value = "* All models have heterogeneous variance"

# ---- Column widths (matches the sibling sheet's bestFit widths) --------
$ws2.Columns.Item(1).ColumnWidth = 16.5
$ws2.Columns.Item(2).ColumnWidth = 9.5
$ws2.Columns.Item(3).ColumnWidth = 2.67
$ws2.Columns.Item(4).ColumnWidth = 12.33
$ws2.Columns.Item(5).ColumnWidth = 16.33
$ws2.Columns.Item(6).ColumnWidth = 16.5
$ws2.Columns.Item(7).ColumnWidth = 20.33

# Final selection / active cell on the new sheet.
$ws2.Range("A7").Select()
